# Update gh-pages output values (F-column "想去人数" counts incremented by 1)
# on sheets "展览" and "全部类型".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 143
$ws1.Range("F10").Value = 1086
$ws1.Range("F15").Value = 388
$ws1.Range("F16").Value = 82
$ws1.Range("F18").Value = 1223

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 143
$ws4.Range("F12").Value = 1086
$ws4.Range("F17").Value = 388
$ws4.Range("F18").Value = 82
$ws4.Range("F20").Value = 1223
